$wb = $excel.ActiveWorkbook

# Updates apply identically to the "展览" and "全部类型" sheets,
# column F ("想去人数"), for rows 3, 4, 6, 8, 9, 10, 12.
$updates = @{
    3  = 1457
    4  = 951
    6  = 2129
    8  = 1290
    9  = 63
    10 = 120
    12 = 309
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
